$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (stored as an Excel date serial number)
$ws.Range("A1").Value = 45436

# Update the price list in column D for rows 33-37
$ws.Range("D33").Value = 9431.802
$ws.Range("D34").Value = 11160.692
$ws.Range("D35").Value = 13785.194
$ws.Range("D36").Value = 15376.297
$ws.Range("D37").Value = 15901.197
